$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.325.21'
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = '1.833.32'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.95%  '
$ws.Range('D5').Value = "'314.73"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('E6').Value = '  +0.81%  '
$ws.Range('D7').Value = "'0.4744"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.80%  '
$ws.Range('D8').Value = "'0.3689"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.93%  '
$ws.Range('D9').Value = "'0.07454"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.28%  '
$ws.Range('D10').Value = "'0.8864"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = '1.881.27'
$ws.Range('E12').Value = '  +4.01%  '
$ws.Range('D13').Value = "'0.07326"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.22%  '
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = "'93.44"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = "'6.579"
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').Value = "'1.011"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = "'0.000008797"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.23%  '
$ws.Range('D19').Value = "'1.011"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.88%  '
$ws.Range('D20').Value = '27.555.40'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').Value = "'5.289"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D23').Value = "'10.67"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.88%  '
$ws.Range('D24').Value = '2.088.01'
$ws.Range('E24').Value = '  +2.01%  '
$ws.Range('D25').Value = "'1.893"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = "'151.94"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('D27').Value = "'18.66"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.87%  '
$ws.Range('D28').Value = "'2.144"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.59%  '
$ws.Range('D29').Value = "'5.236"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = "'117.17"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.31%  '
$ws.Range('D31').Value = "'0.08995"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').Value = "'0.7510"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('D34').Value = "'4.546"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.66%  '
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('D37').Value = "'1.102"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('D38').Value = "'0.05345"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').Value = "'0.01956"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'2.977"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = "'7.266"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  +3.67%  '
$ws.Range('E43').Value = '  +0.66%  '
$ws.Range('E44').Value = '  +0.38%  '
$ws.Range('D45').Value = "'8.484"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('D46').Value = "'0.4915"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.52%  '
$ws.Range('D47').Value = "'10.57"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.26%  '
$ws.Range('D48').Value = "'105.08"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.91%  '
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('D50').Value = "'1.672"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.82%  '
$ws.Range('D51').Value = "'0.06298"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.16%  '
